# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 8 (pushing existing rows 8-24 down to 9-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8; existing rows 8:24 shift down to 9:25.
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the new weekly record.
$ws.Range("A8").Value = 11
$ws.Range("B8").Value = 'Vega Monumental Concepción'
$ws.Range("C8").Value = 'Bíobío'
$ws.Range("D8").Value = 45076
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 'Fruta'
$ws.Range("G8").Value = 100104
$ws.Range("H8").Value = 'Frutos de pepita'
$ws.Range("I8").Value = 100104003
$ws.Range("J8").Value = 'Membrillo'
$ws.Range("K8").Value = 'Champion'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 150
$ws.Range("N8").Value = 10000
$ws.Range("O8").Value = 11000
$ws.Range("P8").Value = 10467
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("R8").Value = 'Provincia de Curicó'
$ws.Range("S8").Value = 582
$ws.Range("T8").Value = 18
